# fix tagihan spp, tagihan biasa dan lainnya
#
# Two students (Badranaya Aksa Wijaya / 2-1-141021 and
# Made Ngurah Prabha Laksmana / 2-1-141042) are removed from the list
# (they had D = 0, i.e. no "Kode Tagihan" yet), the remaining rows shift
# up, and the "Kode Tagihan" values for several students are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two old rows (141021 - Badranaya Aksa Wijaya, 141042 - Made
# Ngurah Prabha Laksmana), which shifts rows 8-18 up to rows 6-16.
$ws.Range("A6:G7").EntireRow.Delete()

# Rows 2-5 (Aidan Maheswara bumi, Salim Alfa Risqi, Aqila Amanina Tertia,
# Azka Syadaad Emiraldi Novriansyah) keep their data/Kode Tagihan as-is.

# Rewrite the remaining rows (now 6-16) with the corrected data.
$data = @(
    @(6,  "2-1-141050", 141050, "Khairi Sarfras Gazala Ibrahim", "EKS-IQR"),
    @(7,  "2-1-141050", 141050, "Khairi Sarfras Gazala Ibrahim", "EKS-RNG"),
    @(8,  "2-1-141031", 141031, "Devendra Hirotta",               "EKS-IQR"),
    @(9,  "2-1-141045", 141045, "Eiliyah Mandara",                "EKS-IQR"),
    @(10, "2-1-141037", 141037, "Muhammad Shah Jellal",           "EKS-IQR"),
    @(11, "2-1-141039", 141039, "Malika Khaira Furqan",           "EKS-IQR"),
    @(12, "2-1-141039", 141039, "Malika Khaira Furqan",           "EKS-RNG"),
    @(13, "2-1-141041", 141041, "Maheswari Ilona Fitranto",       "EKS-IQR"),
    @(14, "2-1-141043", 141043, "Razaneal Zaviyar Wiranatakusumah","EKS-IQR"),
    @(15, "2-1-141034", 141034, "Razatta Muhammad Kamaquinza",    "EKS-IQR"),
    @(16, "2-1-141053", 141053, "Salim Alfa Risqi",               "EKS-RNG")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Update the saved selection to match the target workbook state.
$ws.Range("C22").Select()
